$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new dictionary rows (gender, job, risk, patient_type_map sections) ---
$ws.Range("A67").Value = 'gender'
$ws.Range("B67").Value = 'เพศ'
$ws.Range("C67").Value = 'ชาย'
$ws.Range("D67").Value = 'Male'
$ws.Range("C68").Value = 'หญิง'
$ws.Range("D68").Value = 'Female'
$ws.Range("A69").Value = 'job'
$ws.Range("B69").Value = 'งาน'
$ws.Range("C69").Value = 'นักท่องเที่ยว'
$ws.Range("D69").Value = 'tourist'
$ws.Range("C70").Value = 'พนักงานบริษัท/โรงงาน'
$ws.Range("D70").Value = 'employee'
$ws.Range("C71").Value = 'พนักงานขับรถโดยสาร/รถตู้/แท๊กซี่'
$ws.Range("D71").Value = 'driver'
$ws.Range("C72").Value = 'พ่อบ้าน/แม่บ้าน/ดูแลบ้าน'
$ws.Range("D72").Value = 'housekeeper'
$ws.Range("C73").Value = 'เด็กเล็ก/ในปกครอง'
$ws.Range("D73").Value = 'kid'
$ws.Range("C74").Value = 'มัคคุเทศก์/ไกด์ทัวร์'
$ws.Range("D74").Value = 'guide'
$ws.Range("C75").Value = 'ค้าขาย/ธุรกิจส่วนตัว'
$ws.Range("D75").Value = 'business'
$ws.Range("C76").Value = 'พยาบาล'
$ws.Range("D76").Value = 'nurse'
$ws.Range("C77").Value = 'ว่างงาน'
$ws.Range("D77").Value = 'unemployed'
$ws.Range("C78").Value = 'นักเรียน/นักศึกษา'
$ws.Range("D78").Value = 'student'
$ws.Range("C79").Value = 'ข้าราชการ/พนักงานของรัฐ / รัฐวิสาหกิจ'
$ws.Range("D79").Value = 'government'
$ws.Range("C80").Value = 'พนักงานในสถานบันเทิง'
$ws.Range("D80").Value = 'entertainment '
$ws.Range("C81").Value = 'HCW (บุคลากรทางการแพทย์)'
$ws.Range("D81").Value = 'HCW'
$ws.Range("C82").Value = 'เจ้าหน้าที่สนามบิน'
$ws.Range("D82").Value = 'airport staff'
$ws.Range("C83").Value = 'รับจ้างทั่วไป / ฟรีแลนซ์'
$ws.Range("D83").Value = 'freelance'
$ws.Range("C84").Value = 'ไม่ระบุ'
$ws.Range("D84").Value = 'unspecified'
$ws.Range("C85").Value = 'อื่นๆ'
$ws.Range("D85").Value = 'other'
$ws.Range("C86").Value = 'พนักงานร้านอาหาร'
$ws.Range("D86").Value = 'restaurant staff'
$ws.Range("C87").Value = 'เจ้าหน้าที่บนเครื่องบิน'
$ws.Range("D87").Value = 'cabin crew'
$ws.Range("C88").Value = 'เกษตรกร (ปลูกพืช)'
$ws.Range("D88").Value = 'Farmer '
$ws.Range("C89").Value = 'พนักงานโรงแรม'
$ws.Range("D89").Value = 'hotel staff'
$ws.Range("C90").Value = 'พนักงานนวด/สปา'
$ws.Range("D90").Value = 'Massage '
$ws.Range("C91").Value = 'พระสงฆ์/สามเณร'
$ws.Range("D91").Value = 'monk'
$ws.Range("C92").Value = 'พนักงานบริษัท'
$ws.Range("D92").Value = 'officer'
$ws.Range("C93").Value = 'อิสระ'
$ws.Range("D93").Value = 'independent staff'
$ws.Range("C94").Value = 'เกษตรกร (เลี้ยงสัตว์)'
$ws.Range("D94").Value = 'raising animals'
$ws.Range("C95").Value = 'กรรมกร'
$ws.Range("D95").Value = 'labor'
$ws.Range("C96").Value = 'พนักงานทำความสะอาด'
$ws.Range("D96").Value = 'cleaning staff'
$ws.Range("C97").Value = 'พยาบาลสูตินารีเวช'
$ws.Range("D97").Value = 'obstetrics nurse'
$ws.Range("C98").Value = 'ประมง/จับสัตว์น้ำ'
$ws.Range("D98").Value = 'Fishing'
$ws.Range("C99").Value = 'แพทย์'
$ws.Range("D99").Value = 'doctor'
$ws.Range("A100").Value = 'risk'
$ws.Range("B100").Value = 'ความเสี่ยง'
$ws.Range("C100").Value = 'Cluster CBI โรงเบียร์ 90'
$ws.Range("D100").Value = 'Cluster CBI beer'
$ws.Range("C101").Value = 'Cluster Memory 90s กทม.'
$ws.Range("D101").Value = 'Cluster Memory 90s bkk'
$ws.Range("C102").Value = 'Cluster New Jazz กทม.'
$ws.Range("D102").Value = 'Cluster New Jazz bkk'
$ws.Range("C103").Value = 'Cluster ตลาดพรพัฒน์'
$ws.Range("D103").Value = 'Cluster pornpat market'
$ws.Range("C104").Value = 'Cluster บางแค'
$ws.Range("D104").Value = 'Cluster bangkae'
$ws.Range("C105").Value = 'ไปสถานที่แออัด เช่น งานแฟร์ คอนเสิร์ต'
$ws.Range("D105").Value = 'crowded place'
$ws.Range("C106").Value = 'ไปสถานที่ชุมชน เช่น ตลาดนัด สถานที่ท่องเที่ยว'
$ws.Range("D106").Value = 'community place'
$ws.Range("C107").Value = 'การค้นหาผู้ป่วยเชิงรุกและค้นหาผู้ติดเชื้อในชุมชน'
$ws.Range("D107").Value = 'proactive search'
$ws.Range("C108").Value = 'คนไทยเดินทางกลับจากต่างประเทศ'
$ws.Range("D108").Value = 'thai return'
$ws.Range("C109").Value = 'คนต่างชาติเดินทางมาจากต่างประเทศ'
$ws.Range("D109").Value = 'foreigners'
$ws.Range("C110").Value = 'ตรวจก่อนทำหัตถการ'
$ws.Range("D110").Value = 'medical procedure'
$ws.Range("C111").Value = 'บุคลากรด้านการแพทย์และสาธารณสุข'
$ws.Range("D111").Value = 'health personnel'
$ws.Range("C112").Value = 'ปอดอักเสบ (Pneumonia)'
$ws.Range("D112").Value = 'Pneumonia'
$ws.Range("C113").Value = 'ผู้ติดเชื้อในประเทศ'
$ws.Range("D113").Value = 'infected in country'
$ws.Range("C114").Value = 'ผู้ที่เดินทางมาจากต่างประเทศ และเข้า AOQ'
$ws.Range("D114").Value = 'AOQ'
$ws.Range("C115").Value = 'ผู้ที่เดินทางมาจากต่างประเทศ และเข้า ASQ/ALQ'
$ws.Range("D115").Value = 'ASQ/ALQ'
$ws.Range("C116").Value = 'ผู้ที่เดินทางมาจากต่างประเทศ และเข้า HQ/AHQ'
$ws.Range("D116").Value = 'HQ/AHQ'
$ws.Range("C117").Value = 'ผู้ที่เดินทางมาจากต่างประเทศ และเข้า OQ'
$ws.Range("D117").Value = 'OQ'
$ws.Range("C118").Value = 'พิธีกรรมทางศาสนา'
$ws.Range("D118").Value = 'religious ritual'
$ws.Range("C119").Value = 'ระบุไม่ได้'
$ws.Range("D119").Value = 'unspecified'
$ws.Range("C120").Value = 'ศูนย์กักกัน ผู้ต้องกัก'
$ws.Range("D120").Value = 'detention center'
$ws.Range("C121").Value = 'สถานบันเทิง'
$ws.Range("D121").Value = 'entertainment place'
$ws.Range("C122").Value = 'สนามมวย'
$ws.Range("D122").Value = 'boxing stadium'
$ws.Range("C123").Value = 'สัมผัสใกล้ชิดกับผู้ป่วยยืนยันรายก่อนหน้านี้'
$ws.Range("D123").Value = 'close contact'
$ws.Range("C124").Value = 'อยู่ระหว่างการสอบสวน'
$ws.Range("D124").Value = 'under investigation'
$ws.Range("C125").Value = 'สัมผัสผู้เดินทางจากต่างประเทศ'
$ws.Range("D125").Value = 'contact with travelers'
$ws.Range("C126").Value = 'อาชีพเสี่ยง เช่น ทำงานในสถานที่แออัด หรือทำงานใกล้ชิดสัมผัสชาวต่างชาติ เป็นต้น'
$ws.Range("D126").Value = 'risky career'
$ws.Range("C127").Value = 'Cluster สมุทรสาคร'
$ws.Range("D127").Value = 'Cluster SKN'
$ws.Range("C129").Value = 'อื่นๆ'
$ws.Range("D129").Value = 'other'
$ws.Range("A130").Value = 'patient_type_map'
$ws.Range("C130").Value = '1.ผู้ป่วย PUI'
$ws.Range("D130").Value = 'PUI'
$ws.Range("C131").Value = '10.อื่นๆ'
$ws.Range("D131").Value = 'other'
$ws.Range("C132").Value = '2.สัมผัสผู้ติดเชื้อ'
$ws.Range("D132").Value = 'contact an infected person'
$ws.Range("C133").Value = '3.ต่างชาติมาจากต่างประเทศ'
$ws.Range("D133").Value = 'foreigners'
$ws.Range("C134").Value = '4.คนไทยมาจากต่างประเทศ'
$ws.Range("D134").Value = 'thai return'
$ws.Range("C135").Value = '5.ลักลอบเข้าประเทศ'
$ws.Range("D135").Value = 'smuggle'
$ws.Range("C136").Value = '6.บุคลากรทางการแพทย์'
$ws.Range("D136").Value = 'health personnel'
$ws.Range("C137").Value = '7.เฝ้าระวัง ARI/pneumonia'
$ws.Range("D137").Value = 'ARI/pneumonia'
$ws.Range("C138").Value = '8.สำรวจกลุ่มเสี่ยง (survey)'
$ws.Range("D138").Value = 'survey'
$ws.Range("C139").Value = '9.ขอตรวจหาเชื้อเอง'
$ws.Range("D139").Value = 'self detect'
$ws.Range("C141").Value = '3.จากต่างประเทศ อยู่ใน Quarantine'
$ws.Range("D141").Value = 'Quarantine'
$ws.Range("C142").Value = '8.ขอตรวจหาเชื้อเอง'
$ws.Range("D142").Value = 'self detect'
$ws.Range("C143").Value = '7.สำรวจกลุ่มเสี่ยง (survey)'
$ws.Range("D143").Value = 'survey'
$ws.Range("C144").Value = '5.บุคลากรทางการแพทย์'
$ws.Range("D144").Value = 'health personnel'
$ws.Range("C145").Value = '6.เฝ้าระวัง ARI/ pneumonia'
$ws.Range("D145").Value = 'ARI/pneumonia'
$ws.Range("C146").Value = '8.ขอตรวจหาเชื้อเอง'
$ws.Range("D146").Value = 'self detect'

# --- Column width adjustments ---
$ws.Columns("C").ColumnWidth = 72.33333333333334
$ws.Columns("D").ColumnWidth = 23.666666666666668

# --- Update selection / view to match final state ---
$ws.Range("C18").Select()
